# Insert a new data row before the current row 173 (pushing rows 173..288
# down to 174..289) and populate the new row with a new Pepino ensalada
# price-report entry for "Limache" origin dated 44596 (2022-02-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 173 — shifts the existing rows 173:288 down
# to 174:289, exactly like the target diff's row renumbering.
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new record's data.
$ws.Range("A173").Value = 3
$ws.Range("B173").Value = 'Femacal de La Calera'
$ws.Range("C173").Value = 'Coquimbo'
$ws.Range("D173").Value = 44596
$ws.Range("E173").Value = 5
$ws.Range("F173").Value = 100112043
$ws.Range("G173").Value = 'Pepino ensalada'
$ws.Range("H173").Value = 'Sin especificar'
$ws.Range("I173").Value = 'Primera'
$ws.Range("J173").Value = 110
$ws.Range("K173").Value = 13000
$ws.Range("L173").Value = 14000
$ws.Range("M173").Value = 13545
$ws.Range("N173").Value = '$/caja 70 unidades'
$ws.Range("O173").Value = 'Limache'
$ws.Range("P173").Value = 194
$ws.Range("Q173").Value = 70
$ws.Range("R173").Value = 'Hortaliza'
